$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new daily rows (2025-10-15, serial 45945) below the existing
# data, one per station, mirroring the layout/styles of the prior day's
# rows (28 = 四方坪站, 29 = 高岭站) so the new rows (30, 31) reuse the same
# cell styles (date / currency / integer formats) instead of creating new
# style entries.
$ws.Range("A28:F28").Copy() | Out-Null
$ws.Range("A30:F30").PasteSpecial(-4122) | Out-Null
$ws.Range("A29:F29").Copy() | Out-Null
$ws.Range("A31:F31").PasteSpecial(-4122) | Out-Null

# Row 30: 四方坪站
$ws.Cells.Item(30, 1).Value = 45945
$ws.Cells.Item(30, 2).Value = "四方坪站"
$ws.Cells.Item(30, 3).Value = 8539.9500000000007
$ws.Cells.Item(30, 4).Value = 7025.9
$ws.Cells.Item(30, 5).Value = 3006.5
$ws.Cells.Item(30, 6).Value = 375

# Row 31: 高岭站
$ws.Cells.Item(31, 1).Value = 45945
$ws.Cells.Item(31, 2).Value = "高岭站"
$ws.Cells.Item(31, 3).Value = 4360.3500000000004
$ws.Cells.Item(31, 4).Value = 3504.27
$ws.Cells.Item(31, 5).Value = 1112.5899999999999
$ws.Cells.Item(31, 6).Value = 152

# The active selection moves to H31 after the new rows are added.
$ws.Range("H31").Select() | Out-Null
